# Feedback_Report.xlsx: add the new feedback row submitted by YATIKA JENA
# (roll number and date are text fields in this sheet, so they must stay
# text and not get auto-converted to a number/date by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Cells.Item($row, 1).Value = "YATIKA JENA"
$ws.Cells.Item($row, 2).Value = "'230102105"
$ws.Cells.Item($row, 3).Value = "Very Poor"
$ws.Cells.Item($row, 4).Value = "Very Poor"
$ws.Cells.Item($row, 5).Value = "Very Poor"
$ws.Cells.Item($row, 6).Value = "hojaaa"
$ws.Cells.Item($row, 7).Value = "'5/17/2025"
$ws.Cells.Item($row, 8).Value = "Unknown"
$ws.Cells.Item($row, 9).Value = "Unknown"
